$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Directory" column (H) values were renamed from the old
# "data/Multimedia_Data/Book/" path to the new lowercase
# "data/multimedia/book/" path for every data row.
$ws.Range("H2:H6").Value = "data/multimedia/book/"

# Reflect the editor's on-screen state at save time: scrolled right so
# column E is the left-most visible column, with H12 as the active cell.
[void]$excel.Goto($ws.Range("E1"), $true)
[void]$ws.Range("H12").Select()
